# Insert a new weekly price record for "Piña" (Pineapple) at Terminal
# Hortofrutícola Agro Chillán. This pushes the existing row 263 (and every
# row after it) down by one, and the freed-up row 263 is filled with the
# new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 263 (and everything below it) down by one row.
$ws.Rows.Item(263).EntireRow.Insert()

# Populate the newly inserted row 263 with the new record.
$ws.Cells.Item(263, 1).Value = 7
$ws.Cells.Item(263, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(263, 3).Value = "Ñuble"
$ws.Cells.Item(263, 4).Value = 45146
$ws.Cells.Item(263, 5).Value = 16
$ws.Cells.Item(263, 6).Value = "Fruta"
$ws.Cells.Item(263, 7).Value = 100108
$ws.Cells.Item(263, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(263, 9).Value = 100108005
$ws.Cells.Item(263, 10).Value = "Piña"
$ws.Cells.Item(263, 11).Value = "Caramelo"
$ws.Cells.Item(263, 12).Value = "Segunda"
$ws.Cells.Item(263, 13).Value = 80
$ws.Cells.Item(263, 14).Value = 20000
$ws.Cells.Item(263, 15).Value = 20000
$ws.Cells.Item(263, 16).Value = 20000
$ws.Cells.Item(263, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(263, 18).Value = "Ecuador"
$ws.Cells.Item(263, 19).Value = 1429
$ws.Cells.Item(263, 20).Value = 14

# Keep the new date cell formatted the same way as the other date cells
# in column D (YYYY-MM-DD HH:MM:SS, same style index as its neighbours).
$ws.Cells.Item(263, 4).NumberFormat = $ws.Cells.Item(264, 4).NumberFormat
